$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.622.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -3.39%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.282.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -5.44%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.06%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''595.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -2.84%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''152.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -9.44%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -0.04%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''3.273.42'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -5.69%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.547'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -8.38%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.172'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -11.00%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''6.88'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -2.74%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.510'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -9.94%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''38.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -12.85%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''0.0000247'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -7.91%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.810.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -5.50%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''67.686.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -3.42%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''3.284.56'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -5.81%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''538.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -8.76%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  -5.20%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''7.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -12.20%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''15.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -12.02%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.765'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -10.80%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''7.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -9.70%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''86.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -10.12%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''13.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -9.99%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''3.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -9.58%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +0.15%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''8.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -5.36%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D30").Value = '''29.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -10.94%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''2.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -3.80%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -6.95%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''6.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -15.28%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''5.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -11.46%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''532.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -6.07%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -0.14%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.0454'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -5.85%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''53.46'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -4.93%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.0864'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -10.29%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''9.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -14.76%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -9.80%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''2.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -11.25%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''2.956.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -9.64%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.271'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -9.35%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.0₃0599'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -14.18%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -7.90%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''27.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -12.19%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''2.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -13.94%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -0.07%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -9.68%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''123.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -6.96%  '
$ws.Range("E51").Style = "Normal"
